$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 ("Impressora 4") -------------------------------------------
# Copy the formatting of row 4 onto row 5 first, so the new row inherits the
# same styles (vertical-center + wrap-text) as the other printer rows.
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = "Impressora 4"
$ws.Range("B5").Value = "10.197.0.11"

# --- Headers (rename "Horário" -> "Horario Inicial", add "Horario Final") -
$ws.Range("C1").Value = "Horario Inicial"

# D1 is a brand new cell; copy C1's formatting first so it picks up the same
# header style (bold/center/wrap) instead of the worksheet default.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Horario Final"

# --- "Horario Inicial" column: every printer starts at 08:40 --------------
$ws.Range("C2").Value = 0.3611111111111111
$ws.Range("C3").Value = 0.3611111111111111
$ws.Range("C4").Value = 0.3611111111111111
$ws.Range("C5").Value = 0.3611111111111111

# --- "Horario Final" column: new, staggered end times ---------------------
$ws.Range("D2").NumberFormat = "h:mm:ss"
$ws.Range("D2").Value = 0.375
$ws.Range("D3").NumberFormat = "h:mm:ss"
$ws.Range("D3").Value = 0.38194444444444442
$ws.Range("D4").NumberFormat = "h:mm:ss"
$ws.Range("D4").Value = 0.38888888888888901
$ws.Range("D5").NumberFormat = "h:mm:ss"
$ws.Range("D5").Value = 0.39583333333333298

# --- Column widths for the two new columns ---------------------------------
$ws.Columns("C").ColumnWidth = 13.7109375
$ws.Columns("D").ColumnWidth = 12.42578125

# --- New blank, wrap-text formatted cell below the table -------------------
$ws.Range("F8").WrapText = $true

# --- Selection moves to C6 --------------------------------------------------
$ws.Range("C6").Select()

Write-Host "edit applied"
